$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '67.125.38'
Set-TextValue 2 5 '  -0.85%  '

Set-TextValue 3 4 '2.469.41'
Set-TextValue 3 5 '  -1.75%  '

Set-TextValue 4 5 '  +0.00%  '

Set-TextValue 5 4 '582.92'
Set-TextValue 5 5 '  -1.41%  '

Set-TextValue 6 4 '168.20'
Set-TextValue 6 5 '  -2.51%  '

Set-TextValue 7 5 '  +0.07%  '

Set-TextValue 8 5 '  -1.91%  '

Set-TextValue 9 4 '2.468.18'
Set-TextValue 9 5 '  -1.73%  '

Set-TextValue 10 5 '  -2.40%  '

Set-TextValue 11 5 '  -0.31%  '

Set-TextValue 12 5 '  -2.65%  '

Set-TextValue 13 5 '  -2.94%  '

Set-TextValue 14 4 '25.50'
Set-TextValue 14 5 '  -3.32%  '

Set-TextValue 15 4 '2.915.51'
Set-TextValue 15 5 '  -0.99%  '

Set-TextValue 16 4 '66.956.25'
Set-TextValue 16 5 '  -0.83%  '

Set-TextValue 17 4 '0.0000170'
Set-TextValue 17 5 '  -4.09%  '

Set-TextValue 18 4 '2.479.43'
Set-TextValue 18 5 '  +0.18%  '

Set-TextValue 19 4 '11.14'
Set-TextValue 19 5 '  -5.76%  '

Set-TextValue 20 4 '7.55'
Set-TextValue 20 5 '  -4.16%  '

Set-TextValue 21 4 '353.78'
Set-TextValue 21 5 '  -3.61%  '

Set-TextValue 22 5 '  -2.83%  '

Set-TextValue 23 5 '  -0.66%  '

Set-TextValue 24 4 '69.16'

Set-TextValue 25 4 '4.23'
Set-TextValue 25 5 '  -7.60%  '

Set-TextValue 26 5 '  -7.18%  '

Set-TextValue 27 4 '9.22'
Set-TextValue 27 5 '  -7.53%  '

Set-TextValue 28 4 '0.997'
Set-TextValue 28 5 '  -0.09%  '

Set-TextValue 29 4 '2.594.37'
Set-TextValue 29 5 '  -0.29%  '

Set-TextValue 30 2 'Bittensor'
Set-TextValue 30 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 30 4 '516.27'
Set-TextValue 30 5 '  -3.45%  '

Set-TextValue 31 2 'PEPE'
Set-TextValue 31 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 31 4 '0.0₃0905'
Set-TextValue 31 5 '  -5.79%  '

Set-TextValue 32 4 '7.73'
Set-TextValue 32 5 '  -7.49%  '

Set-TextValue 33 5 '  -5.01%  '

Set-TextValue 34 5 '  -5.96%  '

Set-TextValue 35 5 '  -0.03%  '

Set-TextValue 37 4 '159.79'
Set-TextValue 37 5 '  +0.79%  '

Set-TextValue 39 5 '  -2.88%  '

Set-TextValue 40 5 '  -5.51%  '

Set-TextValue 41 5 '  -0.13%  '

Set-TextValue 42 5 '  -6.43%  '

Set-TextValue 43 5 '  -6.61%  '

Set-TextValue 44 5 '  -6.52%  '

Set-TextValue 45 5 '  -5.08%  '

Set-TextValue 46 4 '38.67'

Set-TextValue 47 4 '140.80'
Set-TextValue 47 5 '  -3.63%  '

Set-TextValue 48 5 '  -6.33%  '

Set-TextValue 49 5 '  -6.55%  '

Set-TextValue 50 5 '  -7.16%  '

Set-TextValue 51 5 '  -10.49%  '
